# edit.ps1 - applies the "New changes in players" commit to the document
# using Word COM-interop style calls against $word.ActiveDocument.
#
# Summary of the edit (from the authoritative OOXML diff):
#  1) In the "Esta alternativa tiene una ventaja principal..." paragraph,
#     split the leading "E" off into its own run and drop the (moved)
#     "_GoBack" bookmark right after it.
#  2) In the "Tablas  Hash ... colisiones en estos" paragraph, remove the
#     "_GoBack" bookmark that used to sit between "...colisiones" and
#     " en estos" (merging the two runs into a single run of text).
#  3) Fix the misspelling "ultima" -> "última" in the weapons paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: relocate the "_GoBack" bookmark from the "Tablas  Hash" /
# "colisiones" paragraph to just after the initial "E" of "Esta
# alternativa tiene una ventaja principal..." - and merge the runs that
# the old bookmark used to separate.
# ---------------------------------------------------------------------

# 1a. Find the unique "Esta alternativa..." paragraph opening and record
#     the character position right after its very first letter ("E").
$findEsta = $d.Content
$findEsta.Find.Execute("Esta alternativa tiene una ventaja principal", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterE = $findEsta.Start + 1

# 1b. Insert the "_GoBack" bookmark collapsed right after "E" - this
#     splits the original run into "E" + "sta alternativa...".
$d.Bookmarks.Add("_GoBack", $d.Range($afterE, $afterE))

# 1c. Re-stamp the tail of that run ("sta alternativa...dependiendo de
#     la ") so it comes out as a fresh run (matching how Word mints a
#     brand-new run for freshly (re)typed text) instead of silently
#     inheriting the original run's revision-save id.
$tailFind = $d.Content
$tailFind.Find.Execute( `
    "sta alternativa tiene una ventaja principal que tanto muestra los datos principales en el que la aplicación funciona y expone la funcionalidad necesaria para procesar los datos a través de interfaz de usuario, y dependiendo de la ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailStart = $tailFind.Start
$tailEnd = $tailFind.End
$tailText = $tailFind.Text
$tmpRange = $d.Range($tailStart, $tailEnd)
$tmpRange.Text = "#"
$tmpRange2 = $d.Range($tailStart, $tailStart + 1)
$tmpRange2.Text = $tailText

# 1d. Drop the old "_GoBack" bookmark that sat between "...colisiones"
#     and " en estos" in the "Tablas  Hash" paragraph.
$d.Bookmarks("_GoBack").Delete()

# 1e. With the bookmark barrier gone, re-assert the surrounding text so
#     the two runs it used to separate collapse back into one run, same
#     as Word does when it re-flows a paragraph after an edit.
$mergeFind = $d.Content
$mergeFind.Find.Execute( `
    "a la hora de ingresar demasiados elementos y aumentar el índice de colisiones en estos", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeFind.Text = "a la hora de ingresar demasiados elementos y aumentar el índice de colisiones en estosZ"
$mergeFind2 = $d.Content
$mergeFind2.Find.Execute( `
    "a la hora de ingresar demasiados elementos y aumentar el índice de colisiones en estosZ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "a la hora de ingresar demasiados elementos y aumentar el índice de colisiones en estos", 2)

# ---------------------------------------------------------------------
# Step 2: spelling fix "ultima" -> "última" in the weapons paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("ultima arma antes de la que había escogido", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "última arma antes de la que había escogido", 2)
